# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the two rows (for
# 52c54f0c-...md and b7e88faa-...md) move from "handed back" state to
# "ready for handoff" state, timestamps advance, and the zh-cn/de-de
# sheets pick up a new "Error Detail" message noting that the handback
# file used wasn't the latest version.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$readyForHandoff = "Ready for handoff"

$errDetail52c54f0c = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f5a2e6b57cc1c9a94779bd70b1d629ea60acdc3/e2e/52c54f0c-3e0c-4864-94c2-7118823c50dd.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2bbc347857331d99ac4a954da72d538bda2c711/e2e/52c54f0c-3e0c-4864-94c2-7118823c50dd.md."
$errDetailB7e88faa = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9f5a2e6b57cc1c9a94779bd70b1d629ea60acdc3/e2e/b7e88faa-6dbd-4f8c-8fde-a41bee98a83d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2bbc347857331d99ac4a954da72d538bda2c711/e2e/b7e88faa-6dbd-4f8c-8fde-a41bee98a83d.md."

# --- Overview sheet: zh-cn / de-de status columns (E, F) + HO Xliff date (G) ---
$overview.Range("E2").Value = $readyForHandoff
$overview.Range("F2").Value = $readyForHandoff
$overview.Range("E3").Value = $readyForHandoff
$overview.Range("F3").Value = $readyForHandoff

$overview.Range("G2").Value = "2016-09-06 03:37:18"
$overview.Range("G3").Value = "2016-09-06 03:37:18"

# --- zh-cn sheet: Status (C), Latest Handoff Datetime (H), Error Detail (P) ---
$zhcn.Range("C2").Value = $readyForHandoff
$zhcn.Range("C3").Value = $readyForHandoff

$zhcn.Range("H2").Value = "2016-09-06 03:37:04"
$zhcn.Range("H3").Value = "2016-09-06 03:37:04"

$zhcn.Range("P2").Value = $errDetail52c54f0c
$zhcn.Range("P3").Value = $errDetailB7e88faa

# --- de-de sheet: Status (C), Error Detail (P) ---
# (H2/H3 on this sheet share the same underlying text as the Overview's
#  G2/G3 "Latest HO Xliff Generate Date" -- keep them in sync.)
$dede.Range("C2").Value = $readyForHandoff
$dede.Range("C3").Value = $readyForHandoff

$dede.Range("H2").Value = "2016-09-06 03:37:18"
$dede.Range("H3").Value = "2016-09-06 03:37:18"

$dede.Range("P2").Value = $errDetail52c54f0c
$dede.Range("P3").Value = $errDetailB7e88faa

# --- Column width adjustments (report re-laid-out for the new, shorter
#     status text and the new, much longer Error Detail column) ---
# Column width is stored in the workbook as (ColumnWidth + 5/6), rounded to
# the nearest 1/6 of a character -- set the closest reachable value.
# (Use numeric column indices -- Columns.Item("<letter>") isn't resolved.)
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332

$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
